# clean languages in contacts df.
# - Rename "Deeper Clean" sheet to "Deeper Cleaning needs"
# - Make the "Deeper Cleaning needs" sheet the active/selected tab
# - Update the selected cell on that sheet from B13 to C20

$wb = $excel.ActiveWorkbook

$ws2 = $wb.Worksheets.Item("Deeper Clean")
$ws2.Name = "Deeper Cleaning needs"

$ws2.Select()
$ws2.Range("C20").Select()
